# [Jihwan] Fix Save System, Add PostProcessing, Add Camera
# Append 5 new Weapon/Sword rows (itemNum 14-18) to the Entities sheet,
# matching the existing rows 12-15 pattern (type=Weapon, weaponType=Sword).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# itemNum (A), type (B), value (D), price (E), weaponType (F)
$newRows = @(
    @{ Row = 16; ItemNum = 14; Value = 4; Price = 300 },
    @{ Row = 17; ItemNum = 15; Value = 5; Price = 300 },
    @{ Row = 18; ItemNum = 16; Value = 6; Price = 300 },
    @{ Row = 19; ItemNum = 17; Value = 7; Price = 300 },
    @{ Row = 20; ItemNum = 18; Value = 8; Price = 300 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $r.ItemNum
    $ws.Range("B$row").Value2 = "Weapon"
    $ws.Range("D$row").Value2 = $r.Value
    $ws.Range("E$row").Value2 = $r.Price
    $ws.Range("F$row").Value2 = "Sword"
}

# Narrow column F (weaponType) - closest achievable width to the authored 13.25
$ws.Columns.Item(6).ColumnWidth = 12.5

# Scroll/selection state left behind by the editor
$ws.Range("F15:F20").Select()
